# surgeDepthProfile1000.xlsx - "data" sheet:
# Column F (sp.cond.us.cm, specific conductivity) was recorded in the wrong
# units (e.g. 0.271) and is corrected to the intended unit convention
# (e.g. 271) for rows 2-14, i.e. every value x1000.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$condValues = @(271, 280, 279, 279, 277, 309, 309, 306, 356, 400, 420, 435, 447)

$row = 2
foreach ($v in $condValues) {
    $ws.Range("F$row").Value = $v
    $row = $row + 1
}

# Leave the sheet scrolled back to the left with the last-touched cell (F7)
# selected, matching the on-disk view state after the edit.
$ws.Range("F7").Select()
